$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H29").Value = 455.05264
$ws.Range("I29").Value = 471.55554
$ws.Range("J29").Value = 440.2
$ws.Range("K29").Value = 1414.66662
$ws.Range("L29").Value = 1320.6
$ws.Range("M29").Value = -1133.66662
$ws.Range("N29").Value = -1882.6
$ws.Range("H34").Value = 1407.6666
$ws.Range("I34").Value = 1407.6666
$ws.Range("K34").Value = 1407.6666
$ws.Range("M34").Value = -1204.6666
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("H36").Value = 1407.6666
$ws.Range("I36").Value = 1407.6666
$ws.Range("K36").Value = 1407.6666
$ws.Range("M36").Value = -692.6666
$ws.Range("H135").Value = 12198546
$ws.Range("I135").Value = 3189.2258
$ws.Range("J135").Value = 50004150
$ws.Range("K135").Value = 28703.0322
$ws.Range("L135").Value = 450037350
$ws.Range("M135").Value = -26168.0322
$ws.Range("N135").Value = -450042420

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 982179.9399999999
$ws.Range("I2").Value = 1712
$ws.Range("J2").Value = 2452882
$ws.Range("K2").Value = 1712
$ws.Range("L2").Value = 2452882
$ws.Range("M2").Value = -1599
$ws.Range("N2").Value = -2453108
$ws.Range("H32").Value = 10761.274
$ws.Range("I32").Value = 6972.795
$ws.Range("J32").Value = 33492.152
$ws.Range("K32").Value = 6972.795
$ws.Range("L32").Value = 33492.152
$ws.Range("M32").Value = -6685.795
$ws.Range("N32").Value = -34066.152
$ws.Range("H59").Value = 35059
$ws.Range("J59").Value = 35059
$ws.Range("L59").Value = 35059
$ws.Range("N59").Value = -36667
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("H116").Value = 982179.9399999999
$ws.Range("I116").Value = 1712
$ws.Range("J116").Value = 2452882
$ws.Range("K116").Value = 1712
$ws.Range("L116").Value = 2452882
$ws.Range("M116").Value = 582
$ws.Range("N116").Value = -2457470

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 982179.9399999999
$ws.Range("I3").Value = 1712
$ws.Range("J3").Value = 2452882
$ws.Range("K3").Value = 1712
$ws.Range("L3").Value = 2452882
$ws.Range("M3").Value = -1598
$ws.Range("N3").Value = -2453110
$ws.Range("H102").Value = 26852
$ws.Range("I102").Value = 15278
$ws.Range("K102").Value = 15278
$ws.Range("M102").Value = -12033
$ws.Range("H105").Value = 1504.7037
$ws.Range("I105").Value = 1435.9565
$ws.Range("J105").Value = 1900
$ws.Range("K105").Value = 1435.9565
$ws.Range("L105").Value = 1900
$ws.Range("M105").Value = 311.0435
$ws.Range("N105").Value = -5394
$ws.Range("H134").Value = 18283054
$ws.Range("I134").Value = 22345556
$ws.Range("J134").Value = 1800
$ws.Range("K134").Value = 67036668
$ws.Range("L134").Value = 5400
$ws.Range("M134").Value = -67034133
$ws.Range("N134").Value = -10470

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1366.5555
$ws.Range("I31").Value = 968.2368
$ws.Range("J31").Value = 3528.8572
$ws.Range("K31").Value = 968.2368
$ws.Range("L31").Value = 3528.8572
$ws.Range("M31").Value = -673.2368
$ws.Range("N31").Value = -4118.8572
$ws.Range("H34").Value = 1366.5555
$ws.Range("I34").Value = 968.2368
$ws.Range("J34").Value = 3528.8572
$ws.Range("K34").Value = 968.2368
$ws.Range("L34").Value = 3528.8572
$ws.Range("M34").Value = -766.2368
$ws.Range("N34").Value = -3932.8572
$ws.Range("H59").Value = 16476.6
$ws.Range("J59").Value = 16595.75
$ws.Range("L59").Value = 16595.75
$ws.Range("N59").Value = -18885.75

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 1611.7273
$ws.Range("I25").Value = 1211.2858
$ws.Range("J25").Value = 2312.5
$ws.Range("K25").Value = 3633.8574
$ws.Range("L25").Value = 6937.5
$ws.Range("M25").Value = -3464.8574
$ws.Range("N25").Value = -7275.5
$ws.Range("H29").Value = 484.2
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 484.2
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 1452.6
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -2006.6
$ws.Range("H30").Value = 1611.7273
$ws.Range("I30").Value = 1211.2858
$ws.Range("J30").Value = 2312.5
$ws.Range("K30").Value = 3633.8574
$ws.Range("L30").Value = 6937.5
$ws.Range("M30").Value = -3531.8574
$ws.Range("N30").Value = -7141.5
$ws.Range("H35").Value = 3933.3333
$ws.Range("J35").Value = 3933.3333
$ws.Range("L35").Value = 11799.9999
$ws.Range("N35").Value = -12375.9999
$ws.Range("H36").Value = 967
$ws.Range("I36").Value = 760.4
$ws.Range("J36").Value = 2000
$ws.Range("K36").Value = 2281.2
$ws.Range("L36").Value = 6000
$ws.Range("M36").Value = -2112.2
$ws.Range("N36").Value = -6338

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 30000
$ws.Range("J48").Value = 30000
$ws.Range("L48").Value = 30000
$ws.Range("N48").Value = -30970
$ws.Range("H132").Value = 3499.75
$ws.Range("I132").Value = 3000
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 9000
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -6470
$ws.Range("N132").Value = -20057

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H132").Value = 5311.325
$ws.Range("I132").Value = 5703.2812
$ws.Range("J132").Value = 3743.5
$ws.Range("K132").Value = 17109.8436
$ws.Range("L132").Value = 11230.5
$ws.Range("M132").Value = -14579.8436
$ws.Range("N132").Value = -16290.5
$ws.Range("H136").Value = 1642.1875
$ws.Range("I136").Value = 1207.7273
$ws.Range("J136").Value = 2598
$ws.Range("K136").Value = 3623.1819
$ws.Range("L136").Value = 7794
$ws.Range("M136").Value = -1073.1819
$ws.Range("N136").Value = -12894

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 10000
$ws.Range("J30").Value = 10000
$ws.Range("L30").Value = 10000
$ws.Range("N30").Value = -10214
$ws.Range("H136").Value = 10124.1
$ws.Range("I136").Value = 14929.192
$ws.Range("J136").Value = 1200.3572
$ws.Range("K136").Value = 44787.576
$ws.Range("L136").Value = 3601.0716
$ws.Range("M136").Value = -42237.576
$ws.Range("N136").Value = -8701.071599999999
